# Update column F (dSF) values for specific rows to match the re-pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 6
    5  = -4
    7  = -3
    13 = 3
    15 = -4
    19 = 5
    20 = -3
    23 = 8
    24 = -6
    28 = 12
    29 = -4
    32 = -1
    37 = -5
    38 = 1
    40 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
